# Generate Report for Handoff
#
# For the files that just completed handoff generation (rows 7 and 10-14 on
# the language sheets), refresh the handoff timestamps and mark their
# Priority as "ht" (hot-fix / handoff-triggered).

$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhcnSheet     = $wb.Worksheets.Item("zh-cn")
$dedeSheet     = $wb.Worksheets.Item("de-de")

$rows = @(7, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview!G  -> "Latest HO Xliff Generate Date"
    $overviewSheet.Cells.Item($r, 7).Value = "2016-08-26 02:21:37"

    # zh-cn!E -> "Priority", zh-cn!H -> "Latest Handoff Datetime"
    $zhcnSheet.Cells.Item($r, 5).Value = "ht"
    $zhcnSheet.Cells.Item($r, 8).Value = "2016-08-26 02:21:32"

    # de-de!E -> "Priority", de-de!H -> "Latest Handoff Datetime"
    $dedeSheet.Cells.Item($r, 5).Value = "ht"
    $dedeSheet.Cells.Item($r, 8).Value = "2016-08-26 02:21:37"
}
